$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates from the latest cryptos data refresh
$ws.Range("D2").Value = "51.706.35"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "3.078.97"
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "388.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0865"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "3.568.45"
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "3.080.38"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.21%  "
$ws.Range("D19").Value = "51.791.79"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.171"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.295"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.39%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("D49").Value = "2.041.81"
$ws.Range("D50").Value = "3.382.30"
$ws.Range("E50").Value = "  +3.16%  "
$ws.Range("E51").Value = "  +6.77%  "
